$wb = $excel.ActiveWorkbook

# Overview sheet: update Status columns for e0757593 row (row 3) from
# "Ready for handoff" to "Handed back: in sync with en-US"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: Status / Latest Handback DateTime / Error Detail for row 3
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-09-06 02:54:16"
$wsZhCn.Range("P3").Value = ""

# de-de sheet: Status / Latest Handback DateTime / Error Detail for row 3
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-09-06 02:54:23"
$wsDeDe.Range("P3").Value = ""

# Error Detail column is now empty, so auto-fit the column width back down
$wsZhCn.Columns.Item(16).AutoFit() | Out-Null
$wsDeDe.Columns.Item(16).AutoFit() | Out-Null
